$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Split the MeDirect "Confirm number to continue using service" sample in two ---
# Row 28 (id 27) originally had motivation "change" (lookalike source, MeDirect) -
# re-classify it as a "lockout" motivation sample.
$ws.Range("F28").Value = "lockout"

# Add a brand new MeDirect sample (id 37) describing the second half of the
# split: a message about a changed number, with a link to a transcript of the
# scam reported in the news.
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = "msg"

# Copy the date cell format from the row above before setting the value, so
# the new cell reuses the existing date style instead of creating a new one.
$ws.Range("C37").Copy()
$ws.Range("C38").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C38").Value = 44384        # 07/07/2021

$ws.Range("D38").Value = "Media"
$ws.Range("E38").Value = "lookalike"
$ws.Range("F38").Value = "change"
$ws.Range("G38").Value = "en"
$ws.Range("H38").Value = "no"
$ws.Range("I38").Value = "number changed, report if not recognised"
$ws.Range("J38").Value = "MeDirect"
$ws.Range("K38").Value = "https://timesofmalta.com/articles/view/scammers-conned-50000-from-40-people-in-one-day.881863"

# --- Update the view state left over from editing the new row ---
$ws.Range("A39").Select()
